# Updates the "cryptos" price/volume table with refreshed data.
# Numeric-looking Price values are written with a leading apostrophe
# (the standard Excel "treat as text" quote-prefix) so that they stay
# text cells (matching the source data, e.g. "67.935.16" which isn't a
# valid number) instead of being auto-converted to floating point
# numbers by Excel's normal type inference.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.935.16'
$ws.Range('E2').Value = '  +2.63%  '

$ws.Range('D3').Value = '3.296.67'
$ws.Range('E3').Value = '  -0.89%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').Value = '''587.78'
$ws.Range('E5').Value = '  +3.28%  '

$ws.Range('D6').Value = '''182.35'
$ws.Range('E6').Value = '  -1.85%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').Value = '''0.587'
$ws.Range('E8').Value = '  +1.73%  '

$ws.Range('D9').Value = '3.288.78'
$ws.Range('E9').Value = '  -0.96%  '

$ws.Range('E10').Value = '  -0.65%  '

$ws.Range('D11').Value = '''0.578'
$ws.Range('E11').Value = '  +0.43%  '

$ws.Range('D12').Value = '''46.10'
$ws.Range('E12').Value = '  -1.23%  '

$ws.Range('D13').Value = '''0.0000273'
$ws.Range('E13').Value = '  +2.58%  '

$ws.Range('D14').Value = '''660.10'
$ws.Range('E14').Value = '  +8.83%  '

$ws.Range('D15').Value = '3.833.37'
$ws.Range('E15').Value = '  -0.70%  '

$ws.Range('D16').Value = '''8.45'
$ws.Range('E16').Value = '  -0.54%  '

$ws.Range('D17').Value = '68.085.56'
$ws.Range('E17').Value = '  +2.70%  '

$ws.Range('E18').Value = '  +1.11%  '

$ws.Range('D19').Value = '3.305.85'
$ws.Range('E19').Value = '  -0.72%  '

$ws.Range('D20').Value = '''17.56'
$ws.Range('E20').Value = '  -2.13%  '

$ws.Range('D21').Value = '''10.92'
$ws.Range('E21').Value = '  -0.89%  '

$ws.Range('D22').Value = '''0.898'
$ws.Range('E22').Value = '  -0.15%  '

$ws.Range('D23').Value = '''17.77'
$ws.Range('E23').Value = '  -1.95%  '

$ws.Range('E24').Value = '  +2.36%  '

$ws.Range('D25').Value = '''97.11'
$ws.Range('E25').Value = '  -2.97%  '

$ws.Range('D26').Value = '''4.01'
$ws.Range('E26').Value = '  +0.04%  '

$ws.Range('D27').Value = '''2.77'
$ws.Range('E27').Value = '  +1.36%  '

$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').Value = '''5.76'
$ws.Range('E28').Value = '  -1.23%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '''9.51'
$ws.Range('E29').Value = '  -0.50%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '''32.42'
$ws.Range('E30').Value = '  +4.39%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '''8.54'
$ws.Range('E31').Value = '  +0.07%  '

$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '''6.79'
$ws.Range('E32').Value = '  +3.09%  '

$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '''600.74'
$ws.Range('E33').Value = '  +7.97%  '

$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '3.945.87'
$ws.Range('E34').Value = '  +3.01%  '

$ws.Range('B35').Value = 'Cosmos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D35').Value = '''10.94'
$ws.Range('E35').Value = '  +0.15%  '

$ws.Range('E36').Value = '  +0.14%  '

$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').Value = '''3.47'
$ws.Range('E37').Value = '  -7.39%  '

$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = '''0.998'
$ws.Range('E38').Value = '  -0.25%  '

$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '''55.76'
$ws.Range('E39').Value = '  -0.93%  '

$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '''0.133'
$ws.Range('E40').Value = '  +2.70%  '

$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '''3.29'
$ws.Range('E41').Value = '  +3.39%  '

$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '''2.67'
$ws.Range('E42').Value = '  +2.31%  '

$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = '''32.97'
$ws.Range('E43').Value = '  -0.53%  '

$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0691'
$ws.Range('E44').Value = '  +0.46%  '

$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').Value = '''3.38'
$ws.Range('E45').Value = '  +0.76%  '

$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').Value = '''0.335'
$ws.Range('E46').Value = '  -0.10%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '''0.0416'
$ws.Range('E47').Value = '  +1.13%  '

$ws.Range('E48').Value = '  +1.22%  '

$ws.Range('E49').Value = '  +0.50%  '

$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '''1.39'
$ws.Range('E50').Value = '  +9.73%  '

$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = '''2.54'
$ws.Range('E51').Value = '  +0.02%  '
